# Auto-generated script to apply profit-sheet numeric updates
# across the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (ALC): "Unbreak My Heart" [Leve Item ID 7015]
$ws.Range("H19").Value = 1228.7778
$ws.Range("I19").Value = 1917
$ws.Range("K19").Value = 1917
$ws.Range("M19").Value = -1742

# Row 86 (ALC): "Filling in the Blanks" [Leve Item ID 12603]
$ws.Range("H86").Value = 5698.3335
$ws.Range("J86").Value = 8745.5
$ws.Range("L86").Value = 8745.5
$ws.Range("N86").Value = -10991.5

# Row 89 (ALC): "Ink into Antiquity (L)" [Leve Item ID 12603]
$ws.Range("H89").Value = 5698.3335
$ws.Range("J89").Value = 8745.5
$ws.Range("L89").Value = 43727.5
$ws.Range("N89").Value = -54959.5

# Row 100 (ALC): "Asking for a Friend" [Leve Item ID 19906]
$ws.Range("H100").Value = 2965.8333
$ws.Range("J100").Value = 1745
$ws.Range("L100").Value = 1745
$ws.Range("N100").Value = -2827

# Row 106 (ALC): "Making Your Mark" [Leve Item ID 19903]
$ws.Range("H106").Value = 4198.8335
$ws.Range("I106").Value = 4198.8335
$ws.Range("K106").Value = 4198.8335
$ws.Range("M106").Value = -3567.8335

# Row 115 (ALC): "5-bell Energy" [Leve Item ID 27957]
$ws.Range("H115").Value = 900
$ws.Range("I115").Value = 900
$ws.Range("K115").Value = 2700
$ws.Range("M115").Value = -1133

# Row 116 (ALC): "Growing Up" [Leve Item ID 27778]
$ws.Range("H116").Value = 2904.111
$ws.Range("I116").Value = 2804
$ws.Range("K116").Value = 2804
$ws.Range("M116").Value = 638

# Row 131 (ALC): "Mindful Study" [Leve Item ID 36108]
$ws.Range("H131").Value = 1656.6666
$ws.Range("I131").Value = 1656.6666
$ws.Range("K131").Value = 4969.9998
$ws.Range("M131").Value = 70.0002000000004

# Row 138 (ALC): "All-night Crafting" [Leve Item ID 44169]
$ws.Range("H138").Value = 2890
$ws.Range("I138").Value = 1181.1111
$ws.Range("J138").Value = 4812.5
$ws.Range("K138").Value = 3543.3333
$ws.Range("L138").Value = 14437.5
$ws.Range("M138").Value = 1596.6667
$ws.Range("N138").Value = -24717.5

$ws = $wb.Worksheets.Item("ARM")
# Row 102 (ARM): "Smells of Rich Tama-hagane" [Leve Item ID 19945]
$ws.Range("H102").Value = 3047.65
$ws.Range("I102").Value = 1530.2667
$ws.Range("K102").Value = 1530.2667
$ws.Range("M102").Value = 91.7333000000001

# Row 110 (ARM): "Scheduled Maintenance" [Leve Item ID 27708]
$ws.Range("H110").Value = 1667.1666
$ws.Range("I110").Value = 1372.5
$ws.Range("J110").Value = 2256.5
$ws.Range("K110").Value = 1372.5
$ws.Range("L110").Value = 2256.5
$ws.Range("M110").Value = 672.5
$ws.Range("N110").Value = -6346.5

# Row 122 (ARM): "Haste for High Durium" [Leve Item ID 36168]
$ws.Range("H122").Value = 575.125
$ws.Range("I122").Value = 575.125
$ws.Range("K122").Value = 1725.375
$ws.Range("M122").Value = 724.625

$ws = $wb.Worksheets.Item("BSM")
# Row 24 (BSM): "Honest Ballast" [Leve Item ID 2420]
$ws.Range("H24").Value = 1800
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1800
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 1800
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -2270

# Row 26 (BSM): "Unseamly Conditions" [Leve Item ID 19535]
$ws.Range("H26").Value = 18238.572
$ws.Range("I26").Value = 18238.572
$ws.Range("K26").Value = 18238.572
$ws.Range("M26").Value = -17946.572

# Row 28 (BSM): "Hearth Maul" [Leve Item ID 19546]
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# Row 96 (BSM): "Hammer Time" [Leve Item ID 19525]
$ws.Range("H96").Value = 15944.571
$ws.Range("I96").Value = 15944.571
$ws.Range("K96").Value = 15944.571
$ws.Range("M96").Value = -13198.571

# Row 107 (BSM): "The Gold Experience" [Leve Item ID 27706]
$ws.Range("H107").Value = 3973.5
$ws.Range("I107").Value = 2768.2
$ws.Range("K107").Value = 2768.2
$ws.Range("M107").Value = -848.1999999999998

# Row 134 (BSM): "Ruthenium Supremium" [Leve Item ID 43998]
$ws.Range("H134").Value = 2152.6
$ws.Range("I134").Value = 939.6667
$ws.Range("K134").Value = 2819.0001
$ws.Range("M134").Value = -284.0001000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 105 (CRP): "Zelkova, My Love" [Leve Item ID 19928]
$ws.Range("H105").Value = 4990
$ws.Range("I105").Value = 4990
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 4990
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3243
$ws.Range("N105").ClearContents()

# Row 107 (CRP): "Built to Last" [Leve Item ID 27689]
$ws.Range("H107").Value = 358.55554
$ws.Range("I107").Value = 224.44444
$ws.Range("J107").Value = 492.66666
$ws.Range("K107").Value = 224.44444
$ws.Range("L107").Value = 492.66666
$ws.Range("M107").Value = 1695.55556
$ws.Range("N107").Value = -4332.66666

# Row 140 (CRP): "Spear Pressure" [Leve Item ID 42455]
$ws.Range("H140").Value = 57167
$ws.Range("J140").Value = 57167
$ws.Range("L140").Value = 57167
$ws.Range("N140").Value = -67527

$ws = $wb.Worksheets.Item("CUL")
# Row 26 (CUL): "A Grape Idea" [Leve Item ID 4746]
$ws.Range("H26").Value = 225
$ws.Range("I26").Value = 225
$ws.Range("K26").Value = 675
$ws.Range("M26").Value = -387

# Row 75 (CUL): "Breakfast of Champions" [Leve Item ID 12863]
$ws.Range("H75").Value = 1500
$ws.Range("I75").Value = 1500
$ws.Range("K75").Value = 4500
$ws.Range("M75").Value = -3502

# Row 78 (CUL): "Emerald Soup for the Soul (L)" [Leve Item ID 12863]
$ws.Range("H78").Value = 1500
$ws.Range("I78").Value = 1500
$ws.Range("K78").Value = 13500
$ws.Range("M78").Value = -8508

# Row 104 (CUL): "Fits to a Tea" [Leve Item ID 19807]
$ws.Range("H104").Value = 9229.23
$ws.Range("J104").Value = 9770.909
$ws.Range("L104").Value = 29312.727
$ws.Range("N104").Value = -34554.727

# Row 140 (CUL): "Sweet, Sweet Bean Juice" [Leve Item ID 44097]
$ws.Range("H140").Value = 1399.8334
$ws.Range("I140").Value = 679.8
$ws.Range("K140").Value = 2039.4
$ws.Range("M140").Value = 3140.6

$ws = $wb.Worksheets.Item("GSM")
# Row 15 (GSM): "The Tusk at Hand" [Leve Item ID 12018]
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 80 (GSM): "Needs More Prayerbell" [Leve Item ID 12521]
$ws.Range("H80").Value = 3138
$ws.Range("I80").Value = 3163.3333
$ws.Range("J80").Value = 3100
$ws.Range("K80").Value = 3163.3333
$ws.Range("L80").Value = 3100
$ws.Range("M80").Value = -2165.3333
$ws.Range("N80").Value = -5096

# Row 81 (GSM): "The Grander Temple" [Leve Item ID 12018]
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 83 (GSM): "With a Noise That Reaches Heaven (L)" [Leve Item ID 12521]
$ws.Range("H83").Value = 3138
$ws.Range("I83").Value = 3163.3333
$ws.Range("J83").Value = 3100
$ws.Range("K83").Value = 15816.6665
$ws.Range("L83").Value = 15500
$ws.Range("M83").Value = -10824.6665
$ws.Range("N83").Value = -25484

# Row 84 (GSM): "Man with a Dragon Earring (L)" [Leve Item ID 12018]
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 109 (GSM): "You're My Wonderhall" [Leve Item ID 25691]
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 113 (GSM): "Copious Crystal Cannons" [Leve Item ID 27710]
$ws.Range("H113").Value = 6963.625
$ws.Range("I113").Value = 3927.5
$ws.Range("K113").Value = 3927.5
$ws.Range("M113").Value = -1757.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW): "Skin off Their Backs" [Leve Item ID 5277]
$ws.Range("H22").Value = 1800.5
$ws.Range("I22").Value = 1800.5
$ws.Range("K22").Value = 1800.5
$ws.Range("M22").Value = -1505.5

# Row 27 (LTW): "Fire and Hide" [Leve Item ID 5277]
$ws.Range("H27").Value = 1800.5
$ws.Range("I27").Value = 1800.5
$ws.Range("K27").Value = 1800.5
$ws.Range("M27").Value = -1693.5

# Row 40 (LTW): "Best Served Toad" [Leve Item ID 36248]
$ws.Range("H40").Value = 4075
$ws.Range("I40").Value = 4227.857
$ws.Range("J40").Value = 3005
$ws.Range("K40").Value = 4227.857
$ws.Range("L40").Value = 3005
$ws.Range("M40").Value = -4091.857
$ws.Range("N40").Value = -3277

# Row 46 (LTW): "Supply Side Logic" [Leve Item ID 5282]
$ws.Range("H46").Value = 6465
$ws.Range("I46").Value = 2000
$ws.Range("K46").Value = 2000
$ws.Range("M46").Value = -1812

# Row 94 (LTW): "Fitting In" [Leve Item ID 18067]
$ws.Range("H94").Value = 12000
$ws.Range("J94").Value = 12000
$ws.Range("L94").Value = 12000
$ws.Range("N94").Value = -13352

# Row 132 (LTW): "Tenets of Tanning" [Leve Item ID 44058]
$ws.Range("H132").Value = 3984.1428
$ws.Range("I132").Value = 2976.8
$ws.Range("J132").Value = 6502.5
$ws.Range("K132").Value = 8930.400000000001
$ws.Range("L132").Value = 19507.5
$ws.Range("M132").Value = -6400.400000000001
$ws.Range("N132").Value = -24567.5

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (WVR): "A Polished Purchase" [Leve Item ID 36210]
$ws.Range("H126").Value = 5220.76
$ws.Range("I126").Value = 3381.3333
$ws.Range("K126").Value = 10143.9999
$ws.Range("M126").Value = -7673.999899999999

# Row 136 (WVR): "Weaving the Envelope" [Leve Item ID 44031]
$ws.Range("H136").Value = 3389.7307
$ws.Range("J136").Value = 4849.5835
$ws.Range("L136").Value = 14548.7505
$ws.Range("N136").Value = -19648.7505
